$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers: F1 = "height", G1 = "weight" (same header style as existing headers)
$ws.Range("F1").Value = "height"
$ws.Range("G1").Value = "weight"
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift the existing "fantasy points" data (column E) out to the new column G,
# then fill the new height/weight columns for every data row.
for ($r = 2; $r -le 17; $r++) {
    $fantasyPoints = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 7).Value = $fantasyPoints
    $ws.Cells.Item($r, 5).Value = 6.25
    $ws.Cells.Item($r, 6).Value = 260
}
